$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cluster Name"
$ws.Range("B1").Value = "Active cases"
$ws.Range("A2").Value = "202111 45784 Holy Rosary Primary SchoolWhite Hills"
$ws.Range("B2").Value = 33
$ws.Range("A3").Value = "3321 Rochester and Elmore District HealthService Yalunkan Aged Care Hostel Rochester"
$ws.Range("B3").Value = 13
$ws.Range("A4").Value = "3323 Villa Maria Catholic Homes St Bernadette'sAged Care Sunshine North"
$ws.Range("B4").Value = 13
$ws.Range("A5").Value = "3600 Belvedere Aged Care Noble Park"
$ws.Range("B5").Value = 21
$ws.Range("A6").Value = "3601 Baptcare Westhaven community outbreak"
$ws.Range("B6").Value = 27
$ws.Range("A7").Value = "3653 Fronditha Thalpori St Albans Aged Care"
$ws.Range("B7").Value = 28
$ws.Range("A8").Value = "4257 BlueCross The Gables Camberwell"
$ws.Range("B8").Value = 27
$ws.Range("A9").Value = "44098 Stawell Primary School"
$ws.Range("B9").Value = 19
$ws.Range("A10").Value = "44121 Wallan Primary School Wallan"
$ws.Range("B10").Value = 22
$ws.Range("A11").Value = "44165 Greenvale Primary School"
$ws.Range("B11").Value = 32
$ws.Range("A12").Value = "44234 Lucknow Primary School Bairnsdale"
$ws.Range("B12").Value = 39
$ws.Range("A13").Value = "44495 Lakes Entrance Primary School"
$ws.Range("B13").Value = 11
$ws.Range("A14").Value = "44667 Beaumaris Primary School Beaumaris"
$ws.Range("B14").Value = 22
$ws.Range("A15").Value = "44701 Hampton Park Primary School HamptonPark"
$ws.Range("B15").Value = 10
$ws.Range("A16").Value = "44718 Parkdale Primary School Parkdale"
$ws.Range("B16").Value = 12
$ws.Range("A17").Value = "44811 Dandenong North Primary SchoolDandenong"
$ws.Range("B17").Value = 32
$ws.Range("A18").Value = "44853 St Albans North Primary School"
$ws.Range("B18").Value = 12
$ws.Range("A19").Value = "44865 Parktone Primary School Parkdale"
$ws.Range("B19").Value = 27
$ws.Range("A20").Value = "44891 Cranbourne Park Primary SchoolCranbourne"
$ws.Range("B20").Value = 19
$ws.Range("A21").Value = "45158 Rowellyn Park Primary School CarrumDowns"
$ws.Range("B21").Value = 11
$ws.Range("A22").Value = "45248 Brookside P-9 College Caroline Springs"
$ws.Range("B22").Value = 13
$ws.Range("A23").Value = "45249 Creekside K-9 College Caroline Springs"
$ws.Range("B23").Value = 17
$ws.Range("A24").Value = "45569 Nhill College Nhill Outbreak"
$ws.Range("B24").Value = 30
$ws.Range("A25").Value = "45648 St Brendans Primary School Shepparton"
$ws.Range("B25").Value = 11
$ws.Range("A26").Value = "4574 Village Glen Aged Care ResidencesMornington"
$ws.Range("B26").Value = 15
$ws.Range("A27").Value = "45967 St Clement of Rome School Bulleen"
$ws.Range("B27").Value = 10
$ws.Range("A28").Value = "46037 Nazareth Catholic Primary SchoolGrovedale"
$ws.Range("B28").Value = 34
$ws.Range("A29").Value = "46050 Our Lady's Catholic Primary SchoolCraigieburn"
$ws.Range("B29").Value = 25
$ws.Range("A30").Value = "46125 Our Lady of the Southern Cross PrimarySchool Manor Lakes"
$ws.Range("B30").Value = 28
$ws.Range("A31").Value = "46190 Haileybury Brighton East"
$ws.Range("B31").Value = 12
$ws.Range("A32").Value = "46215 Yeshivah Primary College St Kilda East"
$ws.Range("B32").Value = 12
$ws.Range("A33").Value = "46276 Hillcrest Christian College Clyde North"
$ws.Range("B33").Value = 18
$ws.Range("A34").Value = "46328 Ilim College Dallas Main Campus DallasOct"
$ws.Range("B34").Value = 27
$ws.Range("A35").Value = "46376 Yesodei HaTorah College Elwood"
$ws.Range("B35").Value = 10
$ws.Range("A36").Value = "50395 St Francis of Assisi Catholic PrimarySchool Tarneit"
$ws.Range("B36").Value = 11
$ws.Range("A37").Value = "52380 Al Iman College Melton South"
$ws.Range("B37").Value = 14
$ws.Range("A38").Value = "52473 John Henry Primary School Pakenham"
$ws.Range("B38").Value = 18
$ws.Range("A39").Value = "Alfred Health Caulfield Hospital Caulfield"
$ws.Range("B39").Value = 14
$ws.Range("A40").Value = "Camp Coolamatong Farm Camp BanksiaPeninsula"
$ws.Range("B40").Value = 12
$ws.Range("A41").Value = "Escala NewQuay Construction Site DocklandsDrive Docklands"
$ws.Range("B41").Value = 17
$ws.Range("A42").Value = "Gippsland and East Gippsland Aboriginal Co-Operative Bairnsdale"
$ws.Range("B42").Value = 15
$ws.Range("A43").Value = "Hamilton Country Music Festival Hamilton GolfClub Hamilton"
$ws.Range("B43").Value = 31
$ws.Range("A44").Value = "Melton Willows Melton"
$ws.Range("B44").Value = 11
$ws.Range("A45").Value = "St Josephs Catholic Primary School Warragul"
$ws.Range("B45").Value = 10
